# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows 64-65) for "Chirimoya" ahead of
# the existing data, shifting the previous rows 64-85 down to rows 66-87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 64 (this pushes old rows 64-85 to 66-87)
$ws.Range("A64:A65").EntireRow.Insert()

# --- New row 64 ---
$ws.Cells.Item(64, 1).Value = 2
$ws.Cells.Item(64, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44798
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100107
$ws.Cells.Item(64, 8).Value = "Otros"
$ws.Cells.Item(64, 9).Value = 100107002
$ws.Cells.Item(64, 10).Value = "Chirimoya"
$ws.Cells.Item(64, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 160
$ws.Cells.Item(64, 14).Value = 23000
$ws.Cells.Item(64, 15).Value = 24000
$ws.Cells.Item(64, 16).Value = 23500
$ws.Cells.Item(64, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(64, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 19).Value = 2350
$ws.Cells.Item(64, 20).Value = 10

# --- New row 65 ---
$ws.Cells.Item(65, 1).Value = 2
$ws.Cells.Item(65, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(65, 3).Value = "Coquimbo"
$ws.Cells.Item(65, 4).Value = 44798
$ws.Cells.Item(65, 5).Value = 4
$ws.Cells.Item(65, 6).Value = "Fruta"
$ws.Cells.Item(65, 7).Value = 100107
$ws.Cells.Item(65, 8).Value = "Otros"
$ws.Cells.Item(65, 9).Value = 100107002
$ws.Cells.Item(65, 10).Value = "Chirimoya"
$ws.Cells.Item(65, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(65, 12).Value = "Segunda"
$ws.Cells.Item(65, 13).Value = 160
$ws.Cells.Item(65, 14).Value = 19000
$ws.Cells.Item(65, 15).Value = 20000
$ws.Cells.Item(65, 16).Value = 19500
$ws.Cells.Item(65, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(65, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(65, 19).Value = 1950
$ws.Cells.Item(65, 20).Value = 10

# Give the inserted date cells the same date format as the rest of column D
$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
